$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values for row 10 (FastPin4 / ATmega2560 build results)
$ws.Range("G10").Value = 360
$ws.Range("H10").Value = 0

# Update the active selection to G10
$ws.Range("G10").Select()
